# Add cell coloring feature: rename sheet, add two new address rows,
# restyle the header (bold) and the Postal code / Country Code columns
# (centered; postal codes kept as text with a quote-prefix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the sheet -------------------------------------------------
$ws.Name = "Sample Addresses"

# --- 2. Insert two new data rows -----------------------------------------
# Row 4 becomes a new "Bad Street" entry; what used to be rows 4-6
# (Beach House, Disney, Poslka Rozana) shift down to 5-6 first, then we
# insert a second new row at what will become row 7 (Brant Home), pushing
# Poslka Rozana down to row 8.
$ws.Rows(4).Insert()
$ws.Rows(7).Insert()

# --- 3. Header row text tweaks --------------------------------------------
$ws.Range("C1").Value = "Postal Code"
$ws.Range("G1").Value = "Notes"

# --- 4. New row 4 data: "Bad Street" entry --------------------------------
$ws.Range("A4").Value = "Bad, so bad"
$ws.Range("B4").Value = "Bad Street"
$ws.Range("C4").Value = "'99999"
$ws.Range("D4").Value = "Wrong City"
$ws.Range("E4").Value = "00.12.34.56"
$ws.Range("F4").Value = "info@bad.er"
$ws.Range("H4").Value = "XY"

# --- 5. New row 7 data: "Brant Home" entry --------------------------------
$ws.Range("A7").Value = "Brant Home"
$ws.Range("B7").Value = "AVENUE DU 8 MAI 1945"
$ws.Range("C7").Value = "'24310"
$ws.Range("D7").Value = "Brantome"
$ws.Range("E7").Value = "55.44.33.22.11"
$ws.Range("F7").Value = "contact@brant.com"
$ws.Range("H7").Value = "FR"

# --- 6. Fill in the previously-missing email for Poslka Rozana (row 8) ---
$ws.Range("F8").Value = "info@rozana.pl"

# --- 7. Re-enter the other postal codes with a quote-prefix so they are
#        stored as text (matches how Excel keeps leading zeros, etc.) ----
$ws.Range("C2").Value = "'10005"
$ws.Range("C3").Value = "'8000"
$ws.Range("C5").Value = "'90254"
$ws.Range("C6").Value = "'32830"
$ws.Range("C8").Value = "'00-791"

# --- 8. Formatting: center the Postal Code and Country Code columns ------
$ws.Range("H2:H8").HorizontalAlignment = -4108
$ws.Range("C2:C8").HorizontalAlignment = -4108

# --- 9. Formatting: bold the header row, with Postal Code/Country Code
#        also centered -----------------------------------------------------
$ws.Range("A1:H1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("H1").HorizontalAlignment = -4108

# --- 10. Column widths: let Excel recompute the best-fit widths now that
#         the content has changed -----------------------------------------
$ws.Columns("A:H").AutoFit()
